$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 135, pushing existing row 135 (and below) down to 136.
$ws.Rows.Item(135).Insert()

# Fill in the new row 135 with data (same as the former row 135 except for the
# columns that changed: D, J, K, L, M, O, P).
$ws.Cells.Item(135, 1).Value = 4
$ws.Cells.Item(135, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(135, 3).Value = "Los Lagos"
$ws.Cells.Item(135, 4).Value = 44574
$ws.Cells.Item(135, 5).Value = 10
$ws.Cells.Item(135, 6).Value = 100112045
$ws.Cells.Item(135, 7).Value = "Zapallo"
$ws.Cells.Item(135, 8).Value = "Paine"
$ws.Cells.Item(135, 9).Value = "1a nueva(o)"
$ws.Cells.Item(135, 10).Value = 600
$ws.Cells.Item(135, 11).Value = 500
$ws.Cells.Item(135, 12).Value = 500
$ws.Cells.Item(135, 13).Value = 500
$ws.Cells.Item(135, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(135, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(135, 16).Value = 500
$ws.Cells.Item(135, 17).Value = 1
$ws.Cells.Item(135, 18).Value = "Hortaliza"
